# Update BMF diet PFAS values (ng/gdw) censored results table on Sheet1.
# The workbook contains a single worksheet with a data table in A1:J13.
# Rows 2-12 (one row per PFAS compound) get refreshed numeric values in
# columns B:J (min/median/max_sole, min/median/max_diet, BMF_diet min/median/max).
# Row 10 (PFTeDA) previously had an "Inf" placeholder in J10 (BMF_diet_max)
# which is now replaced by a concrete, computed numeric value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 0.13
    "C2" = 0.15
    "D2" = 0.21
    "E2" = 0.08
    "F2" = 0.12
    "G2" = 0.24
    "H2" = 0.55
    "I2" = 1.21
    "J2" = 2.66
    "B3" = 1.23
    "C3" = 1.87
    "D3" = 2.8
    "E3" = 0.92
    "F3" = 1.22
    "G3" = 1.36
    "H3" = 0.9
    "I3" = 1.53
    "J3" = 3.04
    "B4" = 0.22
    "C4" = 0.29
    "D4" = 0.55
    "E4" = 0.28
    "F4" = 0.43
    "G4" = 0.86
    "H4" = 0.26
    "I4" = 0.68
    "J4" = 1.98
    "B5" = 0.49
    "C5" = 0.71
    "D5" = 1.06
    "E5" = 2.07
    "F5" = 2.53
    "G5" = 3.94
    "H5" = 0.13
    "I5" = 0.28
    "J5" = 0.51
    "E6" = 4.69
    "F6" = 6.8
    "G6" = 9.76
    "H6" = 0.91
    "I6" = 1.56
    "J6" = 2.93
    "B7" = 1.25
    "C7" = 1.61
    "D7" = 2.16
    "E7" = 0.93
    "F7" = 1.15
    "G7" = 1.48
    "H7" = 0.85
    "I7" = 1.4
    "J7" = 2.32
    "E8" = 0.5600000000000001
    "F8" = 0.72
    "G8" = 0.98
    "H8" = 0.67
    "I8" = 1.4
    "J8" = 2.19
    "E9" = 0.45
    "F9" = 0.57
    "G9" = 1.03
    "H9" = 0.87
    "I9" = 1.99
    "J9" = 4.05
    "B10" = 0.51
    "C10" = 0.92
    "D10" = 0.99
    "E10" = 0.63
    "F10" = 0.9
    "G10" = 1.2
    "H10" = 0.42
    "I10" = 1.01
    "J10" = 1.58
    "B11" = 0.31
    "C11" = 0.42
    "D11" = 0.45
    "E11" = 0.19
    "F11" = 0.23
    "G11" = 0.31
    "H11" = 1
    "I11" = 1.8
    "J11" = 2.43
    "B12" = 0.41
    "C12" = 0.42
    "D12" = 0.52
    "E12" = 0.24
    "F12" = 0.31
    "G12" = 0.46
    "H12" = 0.89
    "I12" = 1.36
    "J12" = 2.15
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
